# Fix/adding translations on the "table_specific_translations" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_specific_translations")

# C3 (text.pt) was mistakenly holding the same "HHID: {{data.hh_id}}" string
# as D3. Correct it to the Portuguese translation.
$ws.Range("C3").Value = "Identificação do agregado{{data.hh_id}}"

# D3 (text.sw) needs its own, distinct Swahili translation instead of
# sharing the Portuguese one.
$ws.Range("D3").Value = "Utambulisho wa Kaya {{data.hh_id}}"

# Both cells now hold longer text, so wrap it and grow the row to fit.
$ws.Range("C3:D3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 25

# Leave the cursor where the editor left it.
$ws.Range("C9").Select() | Out-Null
